# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gains three new trailing columns:
#   H: date             -> "2011-11-22" (text, same value on every data row)
#   I: legislator_name  -> "丁守中"      (text, same value on every data row)
#   J: legislator_id    -> 515           (number, same value on every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legDate = "2011-11-22"
$legName = "丁守中"
$legId = 515

$lastRow = 11

# --- Header row (row 1): new bold/bordered header cells, matching the
# --- existing header style used by columns B..G (copy format from G1). ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (2..11): same three values repeated down every row. ---
for ($r = 2; $r -le $lastRow; $r++) {
    # Assigning a "yyyy-mm-dd"-shaped literal straight to .Value lets Excel's
    # autoconvert silently turn it into a date serial number, so round-trip
    # it through a text formula first and then freeze the formula result
    # into a plain value (Copy + PasteSpecial values) to keep it a string.
    $ws.Cells.Item($r, 8).Formula = "=""$legDate"""
    $ws.Cells.Item($r, 8).Copy()
    $ws.Cells.Item($r, 8).PasteSpecial(-4163)

    $ws.Cells.Item($r, 9).Value = $legName
    $ws.Cells.Item($r, 10).Value = $legId
}
$excel.CutCopyMode = $false

Write-Host "added date/legislator_name/legislator_id columns to 股票 sheet"
